$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Copy the formatting + values of row 5 into the new row 6
$ws.Range("A5:T5").Copy() | Out-Null
$ws.Range("A6:T6").PasteSpecial(-4104) | Out-Null  # xlPasteAll

# Set the new value for column A in row 6 (stored as text "5", like A5="4")
$ws.Range("A6").Value = "5"

# Add the hyperlink for the new row's correoUsuario cell (N6), matching the others
$ws.Hyperlinks.Add($ws.Range("N6"), "mailto:jalzate@todo1.net") | Out-Null

# Update the selection to reflect the new active cell after adding the row
$ws.Range("A7").Select() | Out-Null
